# Refresh market-price-derived leve profit columns (H:N) per item/sheet.
# H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
# K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him | Enchanted Mythrite Ink
$ws.Range("H62").Value = 16111.667
$ws.Range("I62").Value = 6668.3335
$ws.Range("J62").Value = 20833.334
$ws.Range("K62").Value = 6668.3335
$ws.Range("L62").Value = 20833.334
$ws.Range("M62").Value = -6044.3335
$ws.Range("N62").Value = -22081.334
# Row 65: Forgery of Convenience (L) | Enchanted Mythrite Ink
$ws.Range("H65").Value = 16111.667
$ws.Range("I65").Value = 6668.3335
$ws.Range("J65").Value = 20833.334
$ws.Range("K65").Value = 33341.6675
$ws.Range("L65").Value = 104166.67
$ws.Range("M65").Value = -30221.6675
$ws.Range("N65").Value = -110406.67
# Row 76: Warding Off Temptation | Enchanted Hardsilver Ink
$ws.Range("H76").Value = 8806.308000000001
$ws.Range("I76").Value = 10748.667
$ws.Range("J76").Value = 7141.4287
$ws.Range("K76").Value = 10748.667
$ws.Range("L76").Value = 7141.4287
$ws.Range("M76").Value = -10433.667
$ws.Range("N76").Value = -7771.4287
# Row 79: The Garden of Arcane Delights (L) | Enchanted Hardsilver Ink
$ws.Range("H79").Value = 8806.308000000001
$ws.Range("I79").Value = 10748.667
$ws.Range("J79").Value = 7141.4287
$ws.Range("K79").Value = 10748.667
$ws.Range("L79").Value = 7141.4287
$ws.Range("M79").Value = -9656.666999999999
$ws.Range("N79").Value = -9325.4287
# Row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Range("H112").Value = 2006.4117
$ws.Range("I112").Value = 1066.3334
$ws.Range("J112").Value = 2207.8572
$ws.Range("K112").Value = 3199.0002
$ws.Range("L112").Value = 6623.571599999999
$ws.Range("M112").Value = -2091.0002
$ws.Range("N112").Value = -8839.571599999999
# Row 115: 5-bell Energy | Competent Craftsman's Syrup
$ws.Range("H115").Value = 639
$ws.Range("I115").Value = 639
$ws.Range("K115").Value = 1917
$ws.Range("M115").Value = -350
# Row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 5818.4165
$ws.Range("I116").Value = 5545.5713
$ws.Range("K116").Value = 5545.5713
$ws.Range("M116").Value = -2103.5713
# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 1263.6364
$ws.Range("I132").Value = 965.5517
$ws.Range("K132").Value = 2896.6551
$ws.Range("M132").Value = -366.6550999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 25007218
$ws.Range("I61").Value = 25006072
$ws.Range("K61").Value = 25006072
$ws.Range("M61").Value = -25005860
# Row 92: Mail It In | High Steel Scale Mail of Fending
$ws.Range("H92").Value = 59333
$ws.Range("J92").Value = 59333
$ws.Range("L92").Value = 59333
$ws.Range("N92").Value = -64325
# Row 96: The Gauntlet Is Cast | High Steel Gauntlets of Fending
$ws.Range("H96").Value = 35000
$ws.Range("J96").Value = 35000
$ws.Range("L96").Value = 35000
$ws.Range("N96").Value = -40492
# Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws.Range("H102").Value = 14459.417
$ws.Range("I102").Value = 11475.25
$ws.Range("J102").Value = 20427.75
$ws.Range("K102").Value = 11475.25
$ws.Range("L102").Value = 20427.75
$ws.Range("M102").Value = -9853.25
$ws.Range("N102").Value = -23671.75
# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 25007218
$ws.Range("I136").Value = 25006072
$ws.Range("K136").Value = 75018216
$ws.Range("M136").Value = -75015666
# Row 139: Backing up My Words | Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""  # clear - column no longer populated

$ws = $wb.Worksheets.Item("BSM")
# Row 100: And My Axe | Doman Iron War Axe
$ws.Range("H100").Value = 34993.332
$ws.Range("J100").Value = 34993.332
$ws.Range("L100").Value = 34993.332
$ws.Range("N100").Value = -37157.332

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 654611.25
$ws.Range("J31").Value = 1168817.1
$ws.Range("L31").Value = 1168817.1
$ws.Range("N31").Value = -1169407.1
# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 654611.25
$ws.Range("J34").Value = 1168817.1
$ws.Range("L34").Value = 1168817.1
$ws.Range("N34").Value = -1169221.1
# Row 55: Ready for a Rematch | Mythril Lance
$ws.Range("H55").Value = 12500
$ws.Range("I55").Value = 15000
$ws.Range("J55").Value = 10000
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 10000
$ws.Range("M55").Value = -14685
$ws.Range("N55").Value = -10630

$ws = $wb.Worksheets.Item("CUL")
# Row 37: I Love Lamprey | Eel Pie
$ws.Range("H37").Value = 69742.5
$ws.Range("J37").Value = 69742.5
$ws.Range("L37").Value = 209227.5
$ws.Range("N37").Value = -209451.5
# Row 39: Bloody Good Tart, This | Blood Currant Tart
$ws.Range("H39").Value = 1874.75
$ws.Range("J39").Value = 2999
$ws.Range("L39").Value = 8997
$ws.Range("N39").Value = -9585
# Row 64: The Aroma of Faith | Baked Onion Soup
$ws.Range("H64").Value = 3820.6667
$ws.Range("I64").Value = 2312
$ws.Range("J64").Value = 4251.7144
$ws.Range("K64").Value = 6936
$ws.Range("L64").Value = 12755.1432
$ws.Range("M64").Value = -6666
$ws.Range("N64").Value = -13295.1432
# Row 67: Soup's On (L) | Baked Onion Soup
$ws.Range("H67").Value = 3820.6667
$ws.Range("I67").Value = 2312
$ws.Range("J67").Value = 4251.7144
$ws.Range("K67").Value = 6936
$ws.Range("L67").Value = 12755.1432
$ws.Range("M67").Value = -6000
$ws.Range("N67").Value = -14627.1432
# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 1611.8823
$ws.Range("I132").Value = 1240.3
$ws.Range("K132").Value = 11162.7
$ws.Range("M132").Value = -8632.699999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 4542.5
$ws.Range("I80").Value = 3986.1428
$ws.Range("J80").Value = 5321.4
$ws.Range("K80").Value = 3986.1428
$ws.Range("L80").Value = 5321.4
$ws.Range("M80").Value = -2988.1428
$ws.Range("N80").Value = -7317.4
# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 4542.5
$ws.Range("I83").Value = 3986.1428
$ws.Range("J83").Value = 5321.4
$ws.Range("K83").Value = 19930.714
$ws.Range("L83").Value = 26607
$ws.Range("M83").Value = -14938.714
$ws.Range("N83").Value = -36591
# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 90921180
$ws.Range("I132").Value = 100003300
$ws.Range("K132").Value = 300009900
$ws.Range("M132").Value = -300007370

$ws = $wb.Worksheets.Item("LTW")
# Row 100: Tiger in the Sack | Tiger Leather
$ws.Range("H100").Value = 4300
$ws.Range("I100").Value = 3800
$ws.Range("J100").Value = 4400
$ws.Range("K100").Value = 3800
$ws.Range("L100").Value = 4400
$ws.Range("M100").Value = -3259
$ws.Range("N100").Value = -5482
# Row 114: A Heady Endeavor | Atrociraptorskin Headgear of Scouting
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""  # clear - column no longer populated
# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 6386.6
$ws.Range("I122").Value = 6140.909
$ws.Range("J122").Value = 7062.25
$ws.Range("K122").Value = 18422.727
$ws.Range("L122").Value = 21186.75
$ws.Range("M122").Value = -15972.727
$ws.Range("N122").Value = -26086.75

$ws = $wb.Worksheets.Item("WVR")
# Row 37: Bet You Anything | Velveteen Sarouel of Gathering
$ws.Range("H37").Value = 15000
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = ""  # clear - column no longer populated
# Row 128: Lightening Up | Scarlet Moko Gaskins of the Rising Dragon
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = ""  # clear - column no longer populated
